$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update Daniel's status note: the tail of the sentence changes from
#    "...in a few days after meeting 3" to "...by meeting 4 - hopefully".
# ---------------------------------------------------------------------------
$oldTail = "in a few days after meeting 3"
$newTail = "by meeting 4 - hopefully"

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute($oldTail, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newTail, 2)

# Locate the paragraph that now holds the updated note so we can work out
# character offsets for the run splits / bookmark we still need to add.
$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -match "^\*Daniel - ") {
        $notePara = $candidate
        break
    }
}

$paraText = $notePara.Range.Text
$paraStart = $notePara.Range.Start
$paraEnd = $notePara.Range.End

$prefix = "*Daniel - "
$middle = "reading A2 feedback - starting work now - "
$newRun = "will have something for review "

$splitAfterPrefix = $paraStart + $prefix.Length
$splitAfterMiddle = $paraStart + $prefix.Length + $middle.Length
$splitAfterNewRun = $paraStart + $prefix.Length + $middle.Length + $newRun.Length

# ---------------------------------------------------------------------------
# 2) Re-establish the "*Daniel - " / "reading A2 feedback..." run boundary
#    that a plain text replace collapses, and split off the freshly typed
#    "will have something for review " / "by meeting 4 - hopefully" runs.
#    Adding (and immediately removing) a throw-away bookmark at a position
#    is enough to force a run break there without leaving any residue.
# ---------------------------------------------------------------------------
$tmp1 = $d.Range($splitAfterPrefix, $splitAfterPrefix)
$d.Bookmarks.Add("TmpSplitA", $tmp1)
$d.Bookmarks.Item("TmpSplitA").Delete()

$tmp2 = $d.Range($splitAfterNewRun, $splitAfterNewRun)
$d.Bookmarks.Add("TmpSplitB", $tmp2)
$d.Bookmarks.Item("TmpSplitB").Delete()

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark (Word's last-edit marker) from its old spot
#    after "Will we need to purchase any?" to the new edit location, which
#    now spans from just before "will have..." through the end of the note
#    paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$paraEndNow = $notePara.Range.End
$goBackRange = $d.Range($splitAfterMiddle, $paraEndNow)
$d.Bookmarks.Add("_GoBack", $goBackRange)
